$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03994766666666667
$ws.Range("H2").Value = 0.119843
$ws.Range("I2").Value = 0.009195128023837375
$ws.Range("J2").Value = 0.009195128023837375
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.155697
$ws.Range("N2").Value = 0.467091
$ws.Range("O2").Value = 0.021288392311201
$ws.Range("P2").Value = 0.021288392311201
$ws.Range("Q2").Value = 0.006219731857
$ws.Range("R2").Value = 0.05597758671300001
$ws.Range("S2").Value = 0.0001957494927231684
$ws.Range("T2").Value = 0.0001957494927231684
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03994766666666667
$ws.Range("H3").Value = 0.119843
$ws.Range("I3").Value = 0.009195128023837375
$ws.Range("J3").Value = 0.009195128023837375
$ws.Range("O3").Value = 0.01134295290047287
$ws.Range("P3").Value = 0.01134295290047287
$ws.Range("Q3").Value = 0.003314018479
$ws.Range("R3").Value = 0.029826166311
$ws.Range("S3").Value = 0.0001042999040882055
$ws.Range("T3").Value = 0.0001042999040882055
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03994766666666667
$ws.Range("H4").Value = 0.119843
$ws.Range("I4").Value = 0.009195128023837375
$ws.Range("J4").Value = 0.009195128023837375
$ws.Range("M4").Value = 7.075048
$ws.Range("N4").Value = 21.225144
$ws.Range("O4").Value = 0.9673686547883261
$ws.Range("P4").Value = 0.9673686547883261
$ws.Range("Q4").Value = 0.2826316591546666
$ws.Range("R4").Value = 2.543684932392
$ws.Range("S4").Value = 0.008895078627026001
$ws.Range("T4").Value = 0.008895078627026001
$ws.Range("I5").Value = 0.6257489364754861
$ws.Range("J5").Value = 0.6257489364754861
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.155697
$ws.Range("N5").Value = 0.467091
$ws.Range("O5").Value = 0.021288392311201
$ws.Range("P5").Value = 0.021288392311201
$ws.Range("Q5").Value = 0.423266602117
$ws.Range("R5").Value = 3.809399419053
$ws.Range("S5").Value = 0.01332118884800694
$ws.Range("T5").Value = 0.01332118884800694
$ws.Range("I6").Value = 0.6257489364754861
$ws.Range("J6").Value = 0.6257489364754861
$ws.Range("O6").Value = 0.01134295290047287
$ws.Range("P6").Value = 0.01134295290047287
$ws.Range("S6").Value = 0.007097840713962427
$ws.Range("T6").Value = 0.007097840713962427
$ws.Range("I7").Value = 0.6257489364754861
$ws.Range("J7").Value = 0.6257489364754861
$ws.Range("M7").Value = 7.075048
$ws.Range("N7").Value = 21.225144
$ws.Range("O7").Value = 0.9673686547883261
$ws.Range("P7").Value = 0.9673686547883261
$ws.Range("Q7").Value = 19.23371373099467
$ws.Range("R7").Value = 173.103423578952
$ws.Range("S7").Value = 0.6053299069135167
$ws.Range("T7").Value = 0.6053299069135167
$ws.Range("G8").Value = 1.585963
$ws.Range("H8").Value = 4.757889
$ws.Range("I8").Value = 0.3650559355006766
$ws.Range("J8").Value = 0.3650559355006766
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.155697
$ws.Range("N8").Value = 0.467091
$ws.Range("O8").Value = 0.021288392311201
$ws.Range("P8").Value = 0.021288392311201
$ws.Range("Q8").Value = 0.246929681211
$ws.Range("R8").Value = 2.222367130899
$ws.Range("S8").Value = 0.007771453970470891
$ws.Range("T8").Value = 0.007771453970470891
$ws.Range("G9").Value = 1.585963
$ws.Range("H9").Value = 4.757889
$ws.Range("I9").Value = 0.3650559355006766
$ws.Range("J9").Value = 0.3650559355006766
$ws.Range("O9").Value = 0.01134295290047287
$ws.Range("P9").Value = 0.01134295290047287
$ws.Range("Q9").Value = 0.131569904517
$ws.Range("R9").Value = 1.184129140653
$ws.Range("S9").Value = 0.004140812282422235
$ws.Range("T9").Value = 0.004140812282422235
$ws.Range("G10").Value = 1.585963
$ws.Range("H10").Value = 4.757889
$ws.Range("I10").Value = 0.3650559355006766
$ws.Range("J10").Value = 0.3650559355006766
$ws.Range("M10").Value = 7.075048
$ws.Range("N10").Value = 21.225144
$ws.Range("O10").Value = 0.9673686547883261
$ws.Range("P10").Value = 0.9673686547883261
$ws.Range("Q10").Value = 11.220764351224
$ws.Range("R10").Value = 100.986879161016
$ws.Range("S10").Value = 0.3531436692477834
$ws.Range("T10").Value = 0.3531436692477834